$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.083.99'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.260.77'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.54%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.88'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.26'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -3.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.406'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.835.40'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '68.053.02'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.93%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.26'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.76%  '

$ws.Range("E16").Value = '  -2.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.252.24'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("E18").Value = '  -2.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.24'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '415.25'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.51'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.20%  '

$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.25'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.65%  '

$ws.Range("E24").Value = '  -2.61%  '

$ws.Range("E25").Value = '  -3.29%  '

$ws.Range("E26").Value = '  -1.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.34'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.60%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.54'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.43'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.81'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.15%  '

$ws.Range("E33").Value = '  -4.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '164.05'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.58%  '

$ws.Range("E35").Value = '  -5.51%  '

$ws.Range("E36").Value = '  -6.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.73'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.793'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.42'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.31'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.630.57'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.41'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0671'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '336.12'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.11'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.70%  '

$ws.Range("E46").Value = '  -3.79%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.20'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.977'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.86%  '

$ws.Range("E49").Value = '  -1.92%  '

$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.45'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.94%  '
